$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a few "NAME" column values to their PascalCase equivalents
# (order matters for shared-string table placement: CurrentOccupation,
# MaritalStatus, then Weight should be appended in that sequence)
$ws.Range("C20").Value = "CurrentOccupation"
$ws.Range("C19").Value = "MaritalStatus"
$ws.Range("C21").Value = "Weight"

# Update the active cell selection on the sheet
$ws.Range("C27").Select()
